$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 7.636417666666667
$ws.Cells.Item(2, 8).Value = 22.909253
$ws.Cells.Item(2, 9).Value = 0.108532481296676
$ws.Cells.Item(2, 10).Value = 0.108532481296676
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 33.380049
$ws.Cells.Item(2, 14).Value = 100.140147
$ws.Cells.Item(2, 15).Value = 0.3891462059670435
$ws.Cells.Item(2, 16).Value = 0.3891462059670435
$ws.Cells.Item(2, 17).Value = 254.903995897799
$ws.Cells.Item(2, 18).Value = 2294.135963080191
$ws.Cells.Item(2, 19).Value = 0.04223500332079057
$ws.Cells.Item(2, 20).Value = 0.04223500332079057

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 7.636417666666667
$ws.Cells.Item(3, 8).Value = 22.909253
$ws.Cells.Item(3, 9).Value = 0.108532481296676
$ws.Cells.Item(3, 10).Value = 0.108532481296676
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 33.85786133333334
$ws.Cells.Item(3, 14).Value = 101.573584
$ws.Cells.Item(3, 15).Value = 0.3947165649764305
$ws.Cells.Item(3, 16).Value = 0.3947165649764305
$ws.Cells.Item(3, 17).Value = 258.552770441417
$ws.Cells.Item(3, 18).Value = 2326.974933972752
$ws.Cells.Item(3, 19).Value = 0.04283956820579263
$ws.Cells.Item(3, 20).Value = 0.04283956820579263

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 7.636417666666667
$ws.Cells.Item(4, 8).Value = 22.909253
$ws.Cells.Item(4, 9).Value = 0.108532481296676
$ws.Cells.Item(4, 10).Value = 0.108532481296676
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 18.53974466666667
$ws.Cells.Item(4, 14).Value = 55.61923400000001
$ws.Cells.Item(4, 15).Value = 0.2161372290565261
$ws.Cells.Item(4, 16).Value = 0.2161372290565261
$ws.Cells.Item(4, 17).Value = 141.5772337080224
$ws.Cells.Item(4, 18).Value = 1274.195103372202
$ws.Cells.Item(4, 19).Value = 0.02345790977009279
$ws.Cells.Item(4, 20).Value = 0.02345790977009279

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 15.103385
$ws.Cells.Item(5, 8).Value = 45.31015499999999
$ws.Cells.Item(5, 9).Value = 0.214656652056136
$ws.Cells.Item(5, 10).Value = 0.214656652056136
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 33.380049
$ws.Cells.Item(5, 14).Value = 100.140147
$ws.Cells.Item(5, 15).Value = 0.3891462059670435
$ws.Cells.Item(5, 16).Value = 0.3891462059670435
$ws.Cells.Item(5, 17).Value = 504.1517313658649
$ws.Cells.Item(5, 18).Value = 4537.365582292784
$ws.Cells.Item(5, 19).Value = 0.08353282173323308
$ws.Cells.Item(5, 20).Value = 0.08353282173323309

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 15.103385
$ws.Cells.Item(6, 8).Value = 45.31015499999999
$ws.Cells.Item(6, 9).Value = 0.214656652056136
$ws.Cells.Item(6, 10).Value = 0.214656652056136
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 33.85786133333334
$ws.Cells.Item(6, 14).Value = 101.573584
$ws.Cells.Item(6, 15).Value = 0.3947165649764305
$ws.Cells.Item(6, 16).Value = 0.3947165649764305
$ws.Cells.Item(6, 17).Value = 511.3683149939467
$ws.Cells.Item(6, 18).Value = 4602.31483494552
$ws.Cells.Item(6, 19).Value = 0.08472853634893883
$ws.Cells.Item(6, 20).Value = 0.08472853634893883

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 15.103385
$ws.Cells.Item(7, 8).Value = 45.31015499999999
$ws.Cells.Item(7, 9).Value = 0.214656652056136
$ws.Cells.Item(7, 10).Value = 0.214656652056136
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 18.53974466666667
$ws.Cells.Item(7, 14).Value = 55.61923400000001
$ws.Cells.Item(7, 15).Value = 0.2161372290565261
$ws.Cells.Item(7, 16).Value = 0.2161372290565261
$ws.Cells.Item(7, 17).Value = 280.0129015023633
$ws.Cells.Item(7, 18).Value = 2520.11611352127
$ws.Cells.Item(7, 19).Value = 0.04639529397396408
$ws.Cells.Item(7, 20).Value = 0.04639529397396408

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 47.62086333333334
$ws.Cells.Item(8, 8).Value = 142.86259
$ws.Cells.Item(8, 9).Value = 0.676810866647188
$ws.Cells.Item(8, 10).Value = 0.676810866647188
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 33.380049
$ws.Cells.Item(8, 14).Value = 100.140147
$ws.Cells.Item(8, 15).Value = 0.3891462059670435
$ws.Cells.Item(8, 16).Value = 0.3891462059670435
$ws.Cells.Item(8, 17).Value = 1589.58675148897
$ws.Cells.Item(8, 18).Value = 14306.28076340073
$ws.Cells.Item(8, 19).Value = 0.2633783809130198
$ws.Cells.Item(8, 20).Value = 0.2633783809130198

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 47.62086333333334
$ws.Cells.Item(9, 8).Value = 142.86259
$ws.Cells.Item(9, 9).Value = 0.676810866647188
$ws.Cells.Item(9, 10).Value = 0.676810866647188
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 33.85786133333334
$ws.Cells.Item(9, 14).Value = 101.573584
$ws.Cells.Item(9, 15).Value = 0.3947165649764305
$ws.Cells.Item(9, 16).Value = 0.3947165649764305
$ws.Cells.Item(9, 17).Value = 1612.340587313618
$ws.Cells.Item(9, 18).Value = 14511.06528582256
$ws.Cells.Item(9, 19).Value = 0.267148460421699
$ws.Cells.Item(9, 20).Value = 0.267148460421699

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 47.62086333333334
$ws.Cells.Item(10, 8).Value = 142.86259
$ws.Cells.Item(10, 9).Value = 0.676810866647188
$ws.Cells.Item(10, 10).Value = 0.676810866647188
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 18.53974466666667
$ws.Cells.Item(10, 14).Value = 55.61923400000001
$ws.Cells.Item(10, 15).Value = 0.2161372290565261
$ws.Cells.Item(10, 16).Value = 0.2161372290565261
$ws.Cells.Item(10, 17).Value = 882.878647006229
$ws.Cells.Item(10, 18).Value = 7945.907823056062
$ws.Cells.Item(10, 19).Value = 0.1462840253124692
$ws.Cells.Item(10, 20).Value = 0.1462840253124692
